$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = 0.244
$ws.Range("I3").Value = 1.013

$ws.Range("F5").Value = 0.854
$ws.Range("I5").Value = 1.869

$ws.Range("F7").Value = 0.185
$ws.Range("I7").Value = 0.8195

$ws.Range("H7").Select()
